# VIP Group Creation Added
# Rename the placeholder VM names to the real VM names and clear the
# leftover border artifact on the first data row's VMName cell, then
# move the active selection the way the author left it.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = "My-Debian"
$ws.Range("A3").Value = "My-CentOS"
$ws.Range("A4").Value = "My-Ubuntu"
$ws.Range("A5").Value = "My-FreeBSD"
$ws.Range("A6").Value = "My-OpenBSD"
$ws.Range("A7").Value = "My-laggyWin"

# A2 drops its hairline border (the other data rows keep theirs).
$ws.Range("A2").Borders.LineStyle = -4142

# Author left the selection on E8 before saving.
$null = $ws.Range("E8").Select()
